# Refresh the crypto price/volume table (columns D = Price, E = Volume(1h))
# with the latest scraped values. Rows 49/50 (EnergySwap <-> Algorand) also
# swap rank position, so their Coin name / Link / Price / Volume are all
# replaced together.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number -> { column letter -> new cell text }
$rowUpdates = @{
    2 = @{ "D" = "26.846.80"; "E" = "  +0.21%  " }
    3 = @{ "D" = "1.641.63"; "E" = "  -0.21%  " }
    4 = @{ "E" = "  -0.41%  " }
    5 = @{ "D" = "218.15"; "E" = "  +0.60%  " }
    6 = @{ "D" = "0.496"; "E" = "  -0.94%  " }
    7 = @{ "E" = "  -0.37%  " }
    8 = @{ "E" = "  -0.48%  " }
    9 = @{ "E" = "  -1.19%  " }
    10 = @{ "D" = "19.24"; "E" = "  +0.50%  " }
    11 = @{ "E" = "  -0.02%  " }
    12 = @{ "D" = "1.871.58"; "E" = "  -0.05%  " }
    13 = @{ "D" = "1.632.51"; "E" = "  -1.07%  " }
    14 = @{ "E" = "  -0.35%  " }
    15 = @{ "D" = "0.526"; "E" = "  -0.19%  " }
    16 = @{ "D" = "65.22"; "E" = "  +1.02%  " }
    17 = @{ "D" = "26.857.66"; "E" = "  +0.22%  " }
    18 = @{ "D" = "0.0₃0729"; "E" = "  -1.19%  " }
    19 = @{ "D" = "215.42"; "E" = "  +0.79%  " }
    20 = @{ "D" = "1.01"; "E" = "  -0.48%  " }
    21 = @{ "E" = "  -0.24%  " }
    23 = @{ "E" = "  -4.55%  " }
    24 = @{ "D" = "9.20"; "E" = "  -1.64%  " }
    25 = @{ "D" = "147.57"; "E" = "  +1.51%  " }
    26 = @{ "D" = "1.01"; "E" = "  -0.60%  " }
    27 = @{ "E" = "  -0.56%  " }
    28 = @{ "D" = "7.17"; "E" = "  +0.90%  " }
    29 = @{ "D" = "15.69"; "E" = "  +0.19%  " }
    30 = @{ "E" = "  +0.01%  " }
    31 = @{ "E" = "  +1.06%  " }
    32 = @{ "E" = "  +1.26%  " }
    34 = @{ "D" = "1.281.61"; "E" = "  -1.19%  " }
    35 = @{ "E" = "  +0.65%  " }
    36 = @{ "D" = "2.44"; "E" = "  -0.11%  " }
    37 = @{ "E" = "  -0.80%  " }
    38 = @{ "D" = "0.531"; "E" = "  -0.68%  " }
    39 = @{ "D" = "0.821"; "E" = "  -0.47%  " }
    40 = @{ "E" = "  -0.41%  " }
    41 = @{ "D" = "0.807"; "E" = "  -0.45%  " }
    42 = @{ "E" = "  -0.31%  " }
    43 = @{ "D" = "1.783.12" }
    44 = @{ "E" = "  -6.23%  " }
    45 = @{ "D" = "92.50"; "E" = "  +1.01%  " }
    46 = @{ "D" = "61.11"; "E" = "  -1.12%  " }
    47 = @{ "D" = "1.59"; "E" = "  -1.51%  " }
    48 = @{ "E" = "  -1.49%  " }
    49 = @{ "B" = "Algorand"; "C" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; "D" = "0.0970"; "E" = "  -0.62%  " }
    50 = @{ "B" = "EnergySwap"; "C" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; "D" = "7.53"; "E" = "  -1.79%  " }
    51 = @{ "E" = "  -0.42%  " }
}

foreach ($rowNum in $rowUpdates.Keys) {
    $cellsForRow = $rowUpdates[$rowNum]
    foreach ($colLetter in $cellsForRow.Keys) {
        $newValue = $cellsForRow[$colLetter]
        $range = $ws.Range("$colLetter$rowNum")

        # Column D holds prices as text (e.g. "1.871.64" thousand-separated,
        # or "218.14" which Excel would otherwise parse as a number). Force
        # plain-decimal-looking values to stay text, same as the source
        # cells (inline string, no custom style), then drop the temporary
        # "@" format so no stray style survives on the cell.
        $looksNumeric = ($colLetter -eq "D") -and ($newValue -match '^\d+(\.\d+)?$')

        if ($looksNumeric) {
            $range.NumberFormat = "@"
            $range.Value = $newValue
            $range.Style = "Normal"
        } else {
            $range.Value = $newValue
        }
    }
}
